$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 705576.4580000015
$ws.Range("B3").Value = 664859.5580000002
$ws.Range("B4").Value = 2359943.454999996
$ws.Range("B5").Value = 551913.8109999999
$ws.Range("B6").Value = 466825.9979999998
$ws.Range("B7").Value = 2350415.607000005
